{"js": "// 1) \"Did you play any sports last weed?\" is currently split across three\n//    runs (\"Did you play any sports last \", \"weed\", \"?\") with spell-check\n//    proofErr markers wrapped around \"weed\". Collapse it back down to a\n//    single plain run with the same visible text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"Did you play any sports last weed?\";\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === target) {\n    paragraph.insertText(target, \"Replace\");\n    break;\n  }\n}\nawait context.sync();\n\n// 2) Append a new paragraph \"zzzzzzzzzzzzzzzzzzzzz\" right after the last\n//    paragraph in the document body (\"I walk and play sports on the\n//    sand\"), matching the Times New Roman formatting already used there.\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = allParagraphs.items[allParagraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\"zzzzzzzzzzzzzzzzzzzzz\", \"After\");\nnewParagraph.font.name = \"Times New Roman\";\nawait context.sync();\n", "ps1": "# 1) \"Did you play any sports last weed?\" is currently split across three\n#    runs (\"Did you play any sports last \", \"weed\", \"?\") with spell-check\n#    proofErr markers wrapped around \"weed\". Use Find/Replace across the\n#    whole phrase so the match spans all three runs and Word collapses the\n#    result back down into a single plain run with the same visible text.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Did you play any sports last weed?\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Did you play any sports last weed?\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Append a new paragraph \"zzzzzzzzzzzzzzzzzzzzz\" right after the last\n#    paragraph in the document body (\"I walk and play sports on the\n#    sand\"). The new paragraph inherits the Times New Roman formatting\n#    already used by that paragraph.\n$paras = $d.Paragraphs\n$lastParagraph = $paras.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$paras = $d.Paragraphs\n$newParagraph = $paras.Last\n$newParagraph.Range.Text = \"zzzzzzzzzzzzzzzzzzzzz\"\n"}
